$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 18.359437
$ws.Range("H2").Value = 55.078311
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 20.19138166666667
$ws.Range("N2").Value = 60.574145
$ws.Range("O2").Value = 0.2647638724437206
$ws.Range("P2").Value = 0.2647638724437206
$ws.Range("Q2").Value = 370.7023996521217
$ws.Range("R2").Value = 3336.321596869095
$ws.Range("S2").Value = 0.2647638724437206
$ws.Range("T2").Value = 0.2647638724437206

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 18.359437
$ws.Range("H3").Value = 55.078311
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 49.48831433333334
$ws.Range("N3").Value = 148.464943
$ws.Range("O3").Value = 0.6489262577427425
$ws.Range("P3").Value = 0.6489262577427425
$ws.Range("Q3").Value = 908.5775892390303
$ws.Range("R3").Value = 8177.198303151273
$ws.Range("S3").Value = 0.6489262577427425
$ws.Range("T3").Value = 0.6489262577427425

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 18.359437
$ws.Range("H4").Value = 55.078311
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.58215
$ws.Range("N4").Value = 19.74645
$ws.Range("O4").Value = 0.08630986981353689
$ws.Range("P4").Value = 0.0863098698135369
$ws.Range("Q4").Value = 120.84456824955
$ws.Range("R4").Value = 1087.60111424595
$ws.Range("S4").Value = 0.08630986981353689
$ws.Range("T4").Value = 0.0863098698135369
